# ajout script map, commentaires plus precis
#
# 1) Remove the (empty) second slide "Titre 1 / Espace reserve du contenu 2"
#    from the deck - it becomes id 258/rId3 (was 256,257,258 -> now 256,258).
# 2) Refresh the fixed footer date ("Espace reserve de la date") from
#    23/06/2021 to 26/06/2021 on the slide master and every slide layout.

$p = $ppt.ActivePresentation

# ppPlaceholderDate
$ppPlaceholderDate = 16

function Set-DatePlaceholderText {
    # NOTE: positional args only - named binding (-shapes / -newText) drops
    # the live COM reference in this PowerShell host (Count becomes 0).
    param($shapes, $newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Type -eq 14) {
            $phf = $shp.PlaceholderFormat
            if ($phf.Type -eq $ppPlaceholderDate) {
                $shp.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

# --- 1) Delete the now-unwanted slide (index 2) ------------------------
$p.Slides.Item(2).Delete()

# --- 2) Update the fixed footer date on the master & every layout ------
$newDate = "26/06/2021"

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes $newDate

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes $newDate
}
